# Update Name of Algo
# Refresh the imputed RandomForest result values (columns B-D) for the
# rows whose re-run produced slightly different numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.114599999999995
$ws.Range("D4").Value = -6.908099999999997
$ws.Range("C7").Value = -13.04289999999999
$ws.Range("B8").Value = 6.806999999999999
$ws.Range("B10").Value = 5.001000000000001
$ws.Range("D11").Value = -7.460600000000002
$ws.Range("B12").Value = 4.662299999999999
$ws.Range("C14").Value = -13.3398
$ws.Range("D14").Value = -8.085299999999993
$ws.Range("C15").Value = -13.82819999999999
$ws.Range("B18").Value = 7.374299999999993
$ws.Range("C18").Value = -13.28409999999999
$ws.Range("D18").Value = -9.005499999999994
$ws.Range("D19").Value = -9.291899999999991
$ws.Range("C20").Value = -12.2005
$ws.Range("D21").Value = -8.356299999999997
$ws.Range("B25").Value = 6.067900000000001
$ws.Range("D27").Value = -8.632400000000001
$ws.Range("C29").Value = -11.90510000000001
$ws.Range("C30").Value = -13.1291
$ws.Range("C31").Value = -12.78400000000001
$ws.Range("D31").Value = -8.917100000000007
$ws.Range("C35").Value = -11.9555
$ws.Range("B37").Value = 8.904200000000003
$ws.Range("D38").Value = -8.567599999999993
$ws.Range("C40").Value = -13.2609
$ws.Range("D42").Value = -8.745299999999995
$ws.Range("C44").Value = -13.8643
$ws.Range("D44").Value = -7.961099999999998
$ws.Range("D47").Value = -7.496599999999998
$ws.Range("C50").Value = -12.83639999999999
$ws.Range("C54").Value = -13.3282
$ws.Range("B55").Value = 6.292699999999997
$ws.Range("D56").Value = -8.296199999999997
$ws.Range("D58").Value = -8.225199999999997
$ws.Range("D65").Value = -7.949399999999998
$ws.Range("B68").Value = 5.695799999999999
$ws.Range("C68").Value = -11.7966
$ws.Range("D73").Value = -7.880399999999998
$ws.Range("C76").Value = -12.8204
$ws.Range("B77").Value = 8.946800000000007
$ws.Range("B78").Value = 8.903400000000007
$ws.Range("B79").Value = 8.717700000000004
$ws.Range("B80").Value = 9.326500000000003
$ws.Range("B81").Value = 5.331600000000001
$ws.Range("B82").Value = 5.772699999999999
$ws.Range("B84").Value = 5.510000000000002
$ws.Range("C87").Value = -14.10209999999999
$ws.Range("C88").Value = -12.78849999999999
$ws.Range("D90").Value = -7.928299999999998
$ws.Range("C92").Value = -11.9604
$ws.Range("D92").Value = -6.664499999999997
$ws.Range("D94").Value = -6.863699999999998
$ws.Range("D95").Value = -7.730199999999999
$ws.Range("C96").Value = -13.31300000000001
$ws.Range("C98").Value = -11.96259999999999
$ws.Range("B101").Value = 8.8309
$ws.Range("C101").Value = -13.24229999999999
$ws.Range("D101").Value = -7.811099999999998
$ws.Range("B102").Value = 8.424500000000004
$ws.Range("C102").Value = -13.3147
